# BB_Build.xlsx edit: add "Icon" and "ShowName" columns (new G, H) before the
# existing "Desc" column (old G, now shifts to I). Row order / A-F content is
# unchanged; only D's Prefab-path text stays the same too. New G/H columns are
# filled with a short "ShowName"-ish id (derived from the Prefab path) and a
# duplicate of the Chinese display name already in Desc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift the old "Desc" column (G) into the new position (I), before
#        we overwrite G with the new "Icon" column. (Plain `.Value` doesn't
#        round-trip reads in this host - use `.Value2`.) ---
for ($r = 1; $r -le 10; $r++) {
    $old = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 9).Value = $old
}

# --- 2) Header row ---
$ws.Range("G1").Value = "Icon"
$ws.Range("H1").Value = "ShowName"

# --- 3) Data rows: new "Icon" (G) + "ShowName" (H) values ---
$icon = @{
    2  = "Altar_1_1"
    3  = "Arena_1_1"
    4  = "Camp_1_1"
    5  = "GoldMine_1_1"
    6  = "Item_hourse_1_1"
    7  = "League_1_1"
    8  = "MagicHourse_1_1"
    9  = "Tower_1_1"
    10 = "Town_1_1"
}
$showName = @{
    2  = "一级祭坛"
    3  = "一级竞技场"
    4  = "一级兵营"
    5  = "一级金矿"
    6  = "一级道具屋"
    7  = "一级公会"
    8  = "一级魔法屋"
    9  = "一级箭塔"
    10 = "一级大厅"
}

foreach ($r in 2..10) {
    $ws.Cells.Item($r, 7).Value = $icon[$r]
    $ws.Cells.Item($r, 8).Value = $showName[$r]
    # brand-new cells (H, I) default to the "General" style; the rest of the
    # data rows use the workbook's Text ("@") style, so match that.
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 9).NumberFormat = "@"
}

# --- 4) Column widths: columns G:I all end up width 11 (matches the merged
#        col-range the canonical file ends up with). ---
$ws.Range("G1:I1").ColumnWidth = 10.29

# --- 5) Selection moves to H10 after the edits. ---
$ws.Range("H10").Select() | Out-Null
